# Fill in the "Actual Output" (G) column for each test-step row with the
# same value as the "Expected Output" (E) column, so the Pass/FAIL formula
# in column I evaluates to "Pass" for every test case, on every device
# worksheet (Laptop, Tablet - Landscape, Tablet - Protrait, Phone).

$wb = $excel.ActiveWorkbook

$rows = @(6, 7, 8, 9, 10, 11, 12, 20, 21, 22, 23, 24, 25, 26)

foreach ($ws in $wb.Worksheets) {
    foreach ($r in $rows) {
        $expected = $ws.Range("E$r").Value()
        $ws.Range("G$r").Value = $expected
    }
}

# Restore the active-cell selection that Excel leaves behind after the
# edits: the first sheet ends up parked on C2, the other three on G20.
$laptop = $wb.Worksheets.Item("Laptop")
$laptop.Activate()
$laptop.Range("C2").Select()

foreach ($name in @("Tablet - Landscape", "Tablet - Protrait", "Phone")) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Activate()
    $sheet.Range("G20").Select()
}

$laptop.Activate()
